$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 423.14285
$ws.Range("I9").Value = 300.4375
$ws.Range("J9").Value = 815.8
$ws.Range("K9").Value = 300.4375
$ws.Range("L9").Value = 815.8
$ws.Range("M9").Value = -131.4375
$ws.Range("N9").Value = -1153.8
$ws.Range("H17").Value = 4879.6665
$ws.Range("J17").Value = 4879.6665
$ws.Range("L17").Value = 14638.9995
$ws.Range("N17").Value = -14974.9995
$ws.Range("H74").Value = 21302.85
$ws.Range("I74").Value = 22786.555
$ws.Range("K74").Value = 22786.555
$ws.Range("M74").Value = -21850.555
$ws.Range("H77").Value = 21302.85
$ws.Range("I77").Value = 22786.555
$ws.Range("K77").Value = 113932.775
$ws.Range("M77").Value = -109252.775
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H125").Value = 64764.117
$ws.Range("I125").Value = 4574.375
$ws.Range("J125").Value = 118266.11
$ws.Range("K125").Value = 41169.375
$ws.Range("L125").Value = 1064394.99
$ws.Range("M125").Value = -38709.375
$ws.Range("N125").Value = -1069314.99
$ws.Range("H137").Value = 2020.9565
$ws.Range("I137").Value = 1768.2307
$ws.Range("J137").Value = 2349.5
$ws.Range("K137").Value = 5304.6921
$ws.Range("L137").Value = 7048.5
$ws.Range("M137").Value = -2754.6921
$ws.Range("N137").Value = -12148.5
$ws.Range("H138").Value = 2162.2222
$ws.Range("I138").Value = 1984.1364
$ws.Range("J138").Value = 2332.5652
$ws.Range("K138").Value = 5952.4092
$ws.Range("L138").Value = 6997.6956
$ws.Range("M138").Value = -812.4092000000001
$ws.Range("N138").Value = -17277.6956

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9882.367
$ws.Range("I32").Value = 5933.1665
$ws.Range("J32").Value = 20818.615
$ws.Range("K32").Value = 5933.1665
$ws.Range("L32").Value = 20818.615
$ws.Range("M32").Value = -5646.1665
$ws.Range("N32").Value = -21392.615
$ws.Range("H45").Value = 3049.6924
$ws.Range("I45").Value = 2464.8
$ws.Range("K45").Value = 2464.8
$ws.Range("M45").Value = -2087.8
$ws.Range("H61").Value = 4044.4614
$ws.Range("I61").Value = 3554.0588
$ws.Range("J61").Value = 7379.2
$ws.Range("K61").Value = 3554.0588
$ws.Range("L61").Value = 7379.2
$ws.Range("M61").Value = -3342.0588
$ws.Range("N61").Value = -7803.2
$ws.Range("H74").Value = 3240.4285
$ws.Range("I74").Value = 1458.0952
$ws.Range("J74").Value = 5913.9287
$ws.Range("K74").Value = 1458.0952
$ws.Range("L74").Value = 5913.9287
$ws.Range("M74").Value = -584.0952
$ws.Range("N74").Value = -7661.9287
$ws.Range("H77").Value = 3240.4285
$ws.Range("I77").Value = 1458.0952
$ws.Range("J77").Value = 5913.9287
$ws.Range("K77").Value = 7290.476
$ws.Range("L77").Value = 29569.6435
$ws.Range("M77").Value = -2922.476
$ws.Range("N77").Value = -38305.64350000001
$ws.Range("H120").Value = 81000
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 5991.8
$ws.Range("I122").Value = 6970.1665
$ws.Range("J122").Value = 4524.25
$ws.Range("K122").Value = 20910.4995
$ws.Range("L122").Value = 13572.75
$ws.Range("M122").Value = -18460.4995
$ws.Range("N122").Value = -18472.75
$ws.Range("H132").Value = 2899.95
$ws.Range("I132").Value = 2062
$ws.Range("K132").Value = 6186
$ws.Range("M132").Value = -3656
$ws.Range("H136").Value = 4044.4614
$ws.Range("I136").Value = 3554.0588
$ws.Range("J136").Value = 7379.2
$ws.Range("K136").Value = 10662.1764
$ws.Range("L136").Value = 22137.6
$ws.Range("M136").Value = -8112.1764
$ws.Range("N136").Value = -27237.6

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2891.423
$ws.Range("I86").Value = 1181
$ws.Range("K86").Value = 1181
$ws.Range("M86").Value = -58
$ws.Range("H89").Value = 2891.423
$ws.Range("I89").Value = 1181
$ws.Range("K89").Value = 5905
$ws.Range("M89").Value = -289
$ws.Range("H126").Value = 72000
$ws.Range("J126").Value = 72000
$ws.Range("L126").Value = 72000
$ws.Range("N126").Value = -81880
$ws.Range("H134").Value = 3612.3809
$ws.Range("I134").Value = 2782.1052
$ws.Range("K134").Value = 8346.3156
$ws.Range("M134").Value = -5811.3156
$ws.Range("H140").Value = 123999.5
$ws.Range("J140").Value = 123999.5
$ws.Range("L140").Value = 123999.5
$ws.Range("N140").Value = -134359.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6816.9214
$ws.Range("I31").Value = 3379.5757
$ws.Range("J31").Value = 13118.723
$ws.Range("K31").Value = 3379.5757
$ws.Range("L31").Value = 13118.723
$ws.Range("M31").Value = -3084.5757
$ws.Range("N31").Value = -13708.723
$ws.Range("H34").Value = 6816.9214
$ws.Range("I34").Value = 3379.5757
$ws.Range("J34").Value = 13118.723
$ws.Range("K34").Value = 3379.5757
$ws.Range("L34").Value = 13118.723
$ws.Range("M34").Value = -3177.5757
$ws.Range("N34").Value = -13522.723
$ws.Range("H42").Value = 11619.5
$ws.Range("I42").Value = 5989
$ws.Range("J42").Value = 17250
$ws.Range("K42").Value = 5989
$ws.Range("L42").Value = 17250
$ws.Range("M42").Value = -5396
$ws.Range("N42").Value = -18436
$ws.Range("H58").Value = 2684.318
$ws.Range("I58").Value = 2217.4
$ws.Range("J58").Value = 3684.8572
$ws.Range("K58").Value = 2217.4
$ws.Range("L58").Value = 3684.8572
$ws.Range("M58").Value = -2014.4
$ws.Range("N58").Value = -4090.8572
$ws.Range("H105").Value = 2538.2
$ws.Range("I105").Value = 2538.2
$ws.Range("K105").Value = 2538.2
$ws.Range("M105").Value = -791.1999999999998
$ws.Range("H122").Value = 4136.1875
$ws.Range("I122").Value = 3129.7778
$ws.Range("K122").Value = 9389.3334
$ws.Range("M122").Value = -6939.3334
$ws.Range("H136").Value = 2684.318
$ws.Range("I136").Value = 2217.4
$ws.Range("J136").Value = 3684.8572
$ws.Range("K136").Value = 6652.200000000001
$ws.Range("L136").Value = 11054.5716
$ws.Range("M136").Value = -4102.200000000001
$ws.Range("N136").Value = -16154.5716

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1881.5454
$ws.Range("J5").Value = 1256.3334
$ws.Range("L5").Value = 3769.0002
$ws.Range("N5").Value = -3993.0002
$ws.Range("H80").Value = 4499.8335
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 4999.75
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 14999.25
$ws.Range("M80").Value = -9564
$ws.Range("N80").Value = -16871.25
$ws.Range("H83").Value = 4499.8335
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 4999.75
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 44997.75
$ws.Range("M83").Value = -26820
$ws.Range("N83").Value = -54357.75
$ws.Range("H132").Value = 1471.909
$ws.Range("I132").Value = 1422.9524
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 12806.5716
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -10276.5716
$ws.Range("N132").Value = -27560
$ws.Range("H134").Value = 2135.8
$ws.Range("I134").Value = 2135.8
$ws.Range("K134").Value = 6407.400000000001
$ws.Range("M134").Value = -1337.400000000001
$ws.Range("H135").Value = 1881.5454
$ws.Range("J135").Value = 1256.3334
$ws.Range("L135").Value = 11307.0006
$ws.Range("N135").Value = -16377.0006
$ws.Range("H137").Value = 4558.8184
$ws.Range("J137").Value = 5564
$ws.Range("L137").Value = 16692
$ws.Range("N137").Value = -26892

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 279.7143
$ws.Range("I107").Value = 334.5
$ws.Range("J107").Value = 206.66667
$ws.Range("K107").Value = 334.5
$ws.Range("L107").Value = 206.66667
$ws.Range("M107").Value = 1585.5
$ws.Range("N107").Value = -4046.66667
$ws.Range("H122").Value = 7643.2
$ws.Range("I122").Value = 3673.3333
$ws.Range("J122").Value = 10289.777
$ws.Range("K122").Value = 11019.9999
$ws.Range("L122").Value = 30869.331
$ws.Range("M122").Value = -8569.999899999999
$ws.Range("N122").Value = -35769.331
$ws.Range("H132").Value = 5402.4053
$ws.Range("I132").Value = 3787.3635
$ws.Range("J132").Value = 7771.1333
$ws.Range("K132").Value = 11362.0905
$ws.Range("L132").Value = 23313.3999
$ws.Range("M132").Value = -8832.0905
$ws.Range("N132").Value = -28373.3999
$ws.Range("H138").Value = 94542.91
$ws.Range("J138").Value = 94542.91
$ws.Range("L138").Value = 94542.91
$ws.Range("N138").Value = -104822.91

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1028.3334
$ws.Range("I16").Value = 958.9
$ws.Range("K16").Value = 958.9
$ws.Range("M16").Value = -788.9
$ws.Range("H40").Value = 7426.0835
$ws.Range("I40").Value = 4924
$ws.Range("J40").Value = 8260.111000000001
$ws.Range("K40").Value = 4924
$ws.Range("L40").Value = 8260.111000000001
$ws.Range("M40").Value = -4788
$ws.Range("N40").Value = -8532.111000000001
$ws.Range("H132").Value = 8166.476
$ws.Range("I132").Value = 7295.5
$ws.Range("J132").Value = 9908.429
$ws.Range("K132").Value = 21886.5
$ws.Range("L132").Value = 29725.287
$ws.Range("M132").Value = -19356.5
$ws.Range("N132").Value = -34785.287
$ws.Range("H136").Value = 6494.844
$ws.Range("I136").Value = 5253.5283
$ws.Range("J136").Value = 9236.083000000001
$ws.Range("K136").Value = 15760.5849
$ws.Range("L136").Value = 27708.249
$ws.Range("M136").Value = -13210.5849
$ws.Range("N136").Value = -32808.249

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3794.9429
$ws.Range("I126").Value = 3924.2593
$ws.Range("K126").Value = 11772.7779
$ws.Range("M126").Value = -9302.777900000001
$ws.Range("H132").Value = 2739.2642
$ws.Range("I132").Value = 2468.1738
$ws.Range("K132").Value = 7404.5214
$ws.Range("M132").Value = -4874.5214
$ws.Range("H136").Value = 4002.4062
$ws.Range("I136").Value = 2534.2222
$ws.Range("K136").Value = 7602.6666
$ws.Range("M136").Value = -5052.6666
